$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 42-45 : add FR/EN code columns (C/E) and SQL-insert formula
# columns (I/J) for the "Ajouter un Lieu / Details / Modifier /
# Supprimer" block. Row 46 ("Retour") loses its content entirely
# (moved/duplicated elsewhere already at row 66) and disappears.
# ---------------------------------------------------------------------

$ws.Range("C42").Value = "FR"
$ws.Range("E42").Value = "EN"
$ws.Range("I42").Formula = '="INSERT INTO Texte(idTexte, codeTexte, codeLangue, Texte) VALUES(NULL,"""&D42&""","""&C42&""","""&D42&""");"'
$ws.Range("J42:J45").Formula = '="INSERT INTO Texte(idTexte, codeTexte, codeLangue, Texte) VALUES(NULL,"""&D42&""","""&E42&""","""&F42&""");"'

$ws.Range("C43").Value = "FR"
$ws.Range("E43").Value = "EN"
$ws.Range("I43:I46").Formula = '="INSERT INTO Texte(idTexte, codeTexte, codeLangue, Texte) VALUES(NULL,"""&D43&""","""&C43&""","""&D43&""");"'

$ws.Range("C44").Value = "FR"
$ws.Range("E44").Value = "EN"

$ws.Range("C45").Value = "FR"
$ws.Range("E45").Value = "EN"

# Row 46 content removed entirely (including the shared-formula tail
# cell I46 that used to close out the I43:I46 group).
$ws.Range("D46").ClearContents()
$ws.Range("F46").ClearContents()
$ws.Range("I46").ClearContents()

# ---------------------------------------------------------------------
# Row 70 becomes the "ListeCommandes" header row, and two brand new
# rows (71-72) are appended for "Date de Commande" / "Date de
# reception" translations.
# ---------------------------------------------------------------------

$ws.Range("A70").Value = "ListeCommandes"
$ws.Range("C70").Value = "FR"
$ws.Range("D70").Value = "Ajouter une Commande"
$ws.Range("E70").Value = "EN"
$ws.Range("F70").Value = "Add an Order"

$ws.Range("C71").Value = "FR"
$ws.Range("D71").Value = "Date de Commande"
$ws.Range("E71").Value = "EN"

$ws.Range("C72").Value = "FR"
$ws.Range("D72").Value = "Date de reception"
$ws.Range("E72").Value = "EN"

$ws.Range("F71").Value = "Date of Order"
$ws.Range("F72").Value = "Reception Date"

$ws.Range("I70:I72").Formula = '="INSERT INTO Texte(idTexte, codeTexte, codeLangue, Texte) VALUES(NULL,"""&D70&""","""&C70&""","""&D70&""");"'
$ws.Range("J70:J72").Formula = '="INSERT INTO Texte(idTexte, codeTexte, codeLangue, Texte) VALUES(NULL,"""&D70&""","""&E70&""","""&F70&""");"'

# ---------------------------------------------------------------------
# Cosmetic sheet-level tweaks captured by the diff.
# ---------------------------------------------------------------------

$ws.Columns.Item(10).ColumnWidth = 115.140625

$ws.Application.ActiveWindow.ScrollRow = 57
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("I71:I72").Select()
